# "create final Rmd for Ohio" - add an integer stop-id column (to match
# the layout already used on the "Sayfa1" sheet) as column A of "Sheet 1",
# shifting the existing stop/address/lat/long data one column to the right.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet 1")
$ws2 = $wb.Worksheets.Item("Sayfa1")

# Insert a new blank column before column A on "Sheet 1"; this pushes the
# old A:D data to B:E, same as the xml diff (col min/max shift by 1, and
# dimension grows from A1:D46 to A1:E46).
$ws1.Columns.Item(1).Insert()

# Populate the new column A with the GTFS stop_id numbers - these are the
# exact same id values already present in column A of "Sayfa1" (note the
# gaps: some stops, e.g. rows 9, 11, 14, 16, 22, 39, have no id at all).
$ids = @{
    2  = 1;  3  = 2;  4  = 4;  5  = 5;  6  = 6;  7  = 7;  8  = 8;
    10 = 10; 12 = 11; 13 = 12; 15 = 14; 17 = 16; 18 = 17;
    19 = 18; 20 = 19; 21 = 20; 24 = 23; 25 = 24; 26 = 25;
    27 = 26; 28 = 27; 29 = 28; 30 = 29; 31 = 30; 32 = 31;
    33 = 32; 34 = 33; 35 = 34; 36 = 35; 37 = 36; 38 = 37;
    40 = 38; 41 = 39; 42 = 40; 43 = 41; 44 = 42; 45 = 43; 46 = 44
}

foreach ($row in $ids.Keys) {
    $ws1.Cells.Item($row, 1).Value = $ids[$row]
}

# --- View / selection bookkeeping, matching the saved workbook state ---

# "Sayfa1" is no longer the active/selected sheet; its selection becomes
# the id+stop columns (A2:B46), anchored at A2.
$ws2.Activate()
$ws2.Range("A2:B46").Select()

# "Sheet 1" becomes the active tab, with B6 selected.
$ws1.Activate()
$ws1.Range("B6").Select()

Write-Host "Restructured 'Sheet 1' with id column; updated selections."
